$wb = $excel.ActiveWorkbook

# --- SFM sheet: add qualification round data (rows 2-7) ---
$sfm = $wb.Worksheets.Item("SFM")
$sfm.Activate()
$sfm.Range("A2").Value = "Vegard Austrheim Vågen og Henrik Eliassen"
$sfm.Range("B2").Value = 8
$sfm.Range("A3").Value = "Andreas Melheim Hansen og Jørgen Skarsmo"
$sfm.Range("B3").Value = 8
$sfm.Range("A4").Value = "Anders J. Svalestuen og Gabriel Kristiansen"
$sfm.Range("B4").Value = 9
$sfm.Range("A5").Value = "Ole Andre Elvebakk og Georg Kongsvik"
$sfm.Range("B5").Value = 9
$sfm.Range("A6").Value = "Magnus Øslebye og Vegard Tangen"
$sfm.Range("B6").Value = 9
$sfm.Range("A7").Value = "Kasper Støen Nerbøvik og Håvard Idland"
$sfm.Range("B7").Value = 10
$sfm.Range("F6").Select() | Out-Null

# --- SFF sheet: add qualification round data (rows 2-7) ---
$sff = $wb.Worksheets.Item("SFF")
$sff.Activate()
$sff.Range("A2").Value = "Maria Hanssen og Cecilie Rabben"
$sff.Range("B2").Value = 7
$sff.Range("A3").Value = "Beata Wilman og Ingrid Hamnes"
$sff.Range("B3").Value = 7
$sff.Range("A4").Value = "Dawn Stewart og Marie Vik"
$sff.Range("B4").Value = 8
$sff.Range("A5").Value = "Victoria Christensen og Helene Rye Martinsen"
$sff.Range("B5").Value = 9
$sff.Range("A6").Value = "Sara Yuzer og Martine Baalsrud"
$sff.Range("B6").Value = 11
$sff.Range("A7").Value = "Frid Kaspersen og Renate Loraas"
$sff.Range("B7").Value = 13
$sff.Range("C19").Select() | Out-Null

# --- ScoreM sheet: update the selected range ---
$scoreM = $wb.Worksheets.Item("ScoreM")
$scoreM.Activate()
$scoreM.Range("A2:A11").Select() | Out-Null

# --- ScoreF sheet: make it the active/selected tab ---
$scoreF = $wb.Worksheets.Item("ScoreF")
$scoreF.Activate()
$scoreF.Range("C7").Select() | Out-Null
